$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers for rows 2 and 3
$ws.Range("Q2").Value = 528943
$ws.Range("R2").Value = 6229759

$ws.Range("Q3").Value = 528908
$ws.Range("R3").Value = 6229764

# Remove the now-unused Starttid (Z) / Sluttid (AB) time values
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
